# Natmi following Dr Hou advice
# Updates the LR-pairs (Icosl-Cd28) NATMI results table:
#  - recomputed stats for the existing sending/target cluster pairs (rows 2-10)
#  - added the missing M2 sending-cluster rows (rows 11-16)
# Dimension grows from A1:T11 to A1:T16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> M1
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Icosl"
$ws.Range("C2").Value = "Cd28"
$ws.Range("D2").Value = "M1"
$ws.Range("E2").Value = 3.0
$ws.Range("F2").Value = 1.0
$ws.Range("G2").Value = 1.155534666666667
$ws.Range("H2").Value = 3.466604
$ws.Range("I2").Value = 0.08588454688409856
$ws.Range("J2").Value = 0.08588454688409856
$ws.Range("K2").Value = 3.0
$ws.Range("L2").Value = 1.0
$ws.Range("M2").Value = 2.3998
$ws.Range("N2").Value = 7.199400000000001
$ws.Range("O2").Value = 0.3693744204955859
$ws.Range("P2").Value = 0.3693744204955859
$ws.Range("Q2").Value = 2.773052093066668
$ws.Range("R2").Value = 24.9574688376
$ws.Range("S2").Value = 0.03172355473483989
$ws.Range("T2").Value = 0.03172355473483988

# Row 3: ECs -> M2
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Icosl"
$ws.Range("C3").Value = "Cd28"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3.0
$ws.Range("F3").Value = 1.0
$ws.Range("G3").Value = 1.155534666666667
$ws.Range("H3").Value = 3.466604
$ws.Range("I3").Value = 0.08588454688409856
$ws.Range("J3").Value = 0.08588454688409856
$ws.Range("K3").Value = 3.0
$ws.Range("L3").Value = 1.0
$ws.Range("M3").Value = 4.065646999999999
$ws.Range("N3").Value = 12.196941
$ws.Range("O3").Value = 0.6257796502061076
$ws.Range("P3").Value = 0.6257796502061076
$ws.Range("Q3").Value = 4.697996050929333
$ws.Range("R3").Value = 42.281964458364
$ws.Range("S3").Value = 0.05374480170724125
$ws.Range("T3").Value = 0.05374480170724125

# Row 4: ECs -> sCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Icosl"
$ws.Range("C4").Value = "Cd28"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3.0
$ws.Range("F4").Value = 1.0
$ws.Range("G4").Value = 1.155534666666667
$ws.Range("H4").Value = 3.466604
$ws.Range("I4").Value = 0.08588454688409856
$ws.Range("J4").Value = 0.08588454688409856
$ws.Range("K4").Value = 1.0
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.03148366666666667
$ws.Range("N4").Value = 0.094451
$ws.Range("O4").Value = 0.004845929298306607
$ws.Range("P4").Value = 0.004845929298306606
$ws.Range("Q4").Value = 0.03638046826711112
$ws.Range("R4").Value = 0.327424214404
$ws.Range("S4").Value = 0.0004161904420174406
$ws.Range("T4").Value = 0.0004161904420174405

# Row 5: FAPs -> M1
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Icosl"
$ws.Range("C5").Value = "Cd28"
$ws.Range("D5").Value = "M1"
$ws.Range("E5").Value = 3.0
$ws.Range("F5").Value = 1.0
$ws.Range("G5").Value = 0.8102856666666666
$ws.Range("H5").Value = 2.430857
$ws.Range("I5").Value = 0.06022408443105677
$ws.Range("J5").Value = 0.06022408443105677
$ws.Range("K5").Value = 3.0
$ws.Range("L5").Value = 1.0
$ws.Range("M5").Value = 2.3998
$ws.Range("N5").Value = 7.199400000000001
$ws.Range("O5").Value = 0.3693744204955859
$ws.Range("P5").Value = 0.3693744204955859
$ws.Range("Q5").Value = 1.944523542866667
$ws.Range("R5").Value = 17.5007118858
$ws.Range("S5").Value = 0.02224523628659883
$ws.Range("T5").Value = 0.02224523628659883

# Row 6: FAPs -> M2
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Icosl"
$ws.Range("C6").Value = "Cd28"
$ws.Range("D6").Value = "M2"
$ws.Range("E6").Value = 3.0
$ws.Range("F6").Value = 1.0
$ws.Range("G6").Value = 0.8102856666666666
$ws.Range("H6").Value = 2.430857
$ws.Range("I6").Value = 0.06022408443105677
$ws.Range("J6").Value = 0.06022408443105677
$ws.Range("K6").Value = 3.0
$ws.Range("L6").Value = 1.0
$ws.Range("M6").Value = 4.065646999999999
$ws.Range("N6").Value = 12.196941
$ws.Range("O6").Value = 0.6257796502061076
$ws.Range("P6").Value = 0.6257796502061076
$ws.Range("Q6").Value = 3.294335489826333
$ws.Range("R6").Value = 29.64901940843699
$ws.Range("S6").Value = 0.0376870064892498
$ws.Range("T6").Value = 0.0376870064892498

# Row 7: FAPs -> sCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Icosl"
$ws.Range("C7").Value = "Cd28"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3.0
$ws.Range("F7").Value = 1.0
$ws.Range("G7").Value = 0.8102856666666666
$ws.Range("H7").Value = 2.430857
$ws.Range("I7").Value = 0.06022408443105677
$ws.Range("J7").Value = 0.06022408443105677
$ws.Range("K7").Value = 1.0
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.03148366666666667
$ws.Range("N7").Value = 0.094451
$ws.Range("O7").Value = 0.004845929298306607
$ws.Range("P7").Value = 0.004845929298306606
$ws.Range("Q7").Value = 0.02551076383411111
$ws.Range("R7").Value = 0.2295968745069999
$ws.Range("S7").Value = 0.0002918416552081488
$ws.Range("T7").Value = 0.0002918416552081487

# Row 8: M1 -> M1
$ws.Range("A8").Value = "M1"
$ws.Range("B8").Value = "Icosl"
$ws.Range("C8").Value = "Cd28"
$ws.Range("D8").Value = "M1"
$ws.Range("E8").Value = 3.0
$ws.Range("F8").Value = 1.0
$ws.Range("G8").Value = 3.685238666666667
$ws.Range("H8").Value = 11.055716
$ws.Range("I8").Value = 0.2739035549313618
$ws.Range("J8").Value = 0.2739035549313618
$ws.Range("K8").Value = 3.0
$ws.Range("L8").Value = 1.0
$ws.Range("M8").Value = 2.3998
$ws.Range("N8").Value = 7.199400000000001
$ws.Range("O8").Value = 0.3693744204955859
$ws.Range("P8").Value = 0.3693744204955859
$ws.Range("Q8").Value = 8.843835752266669
$ws.Range("R8").Value = 79.59452177040001
$ws.Range("S8").Value = 0.1011729668744527
$ws.Range("T8").Value = 0.1011729668744527

# Row 9: M1 -> M2
$ws.Range("A9").Value = "M1"
$ws.Range("B9").Value = "Icosl"
$ws.Range("C9").Value = "Cd28"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 3.0
$ws.Range("F9").Value = 1.0
$ws.Range("G9").Value = 3.685238666666667
$ws.Range("H9").Value = 11.055716
$ws.Range("I9").Value = 0.2739035549313618
$ws.Range("J9").Value = 0.2739035549313618
$ws.Range("K9").Value = 3.0
$ws.Range("L9").Value = 1.0
$ws.Range("M9").Value = 4.065646999999999
$ws.Range("N9").Value = 12.196941
$ws.Range("O9").Value = 0.6257796502061076
$ws.Range("P9").Value = 0.6257796502061076
$ws.Range("Q9").Value = 14.98287952941733
$ws.Range("R9").Value = 134.845915764756
$ws.Range("S9").Value = 0.171403270795157
$ws.Range("T9").Value = 0.171403270795157

# Row 10: M1 -> sCs
$ws.Range("A10").Value = "M1"
$ws.Range("B10").Value = "Icosl"
$ws.Range("C10").Value = "Cd28"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3.0
$ws.Range("F10").Value = 1.0
$ws.Range("G10").Value = 3.685238666666667
$ws.Range("H10").Value = 11.055716
$ws.Range("I10").Value = 0.2739035549313618
$ws.Range("J10").Value = 0.2739035549313618
$ws.Range("K10").Value = 1.0
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.03148366666666667
$ws.Range("N10").Value = 0.094451
$ws.Range("O10").Value = 0.004845929298306607
$ws.Range("P10").Value = 0.004845929298306606
$ws.Range("Q10").Value = 0.1160248257684444
$ws.Range("R10").Value = 1.044223431916
$ws.Range("S10").Value = 0.001327317261752219
$ws.Range("T10").Value = 0.001327317261752219

# Row 11: M2 -> M1
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Icosl"
$ws.Range("C11").Value = "Cd28"
$ws.Range("D11").Value = "M1"
$ws.Range("E11").Value = 3.0
$ws.Range("F11").Value = 1.0
$ws.Range("G11").Value = 7.390795
$ws.Range("H11").Value = 22.172385
$ws.Range("I11").Value = 0.5493172104644152
$ws.Range("J11").Value = 0.5493172104644152
$ws.Range("K11").Value = 3.0
$ws.Range("L11").Value = 1.0
$ws.Range("M11").Value = 2.3998
$ws.Range("N11").Value = 7.199400000000001
$ws.Range("O11").Value = 0.3693744204955859
$ws.Range("P11").Value = 0.3693744204955859
$ws.Range("Q11").Value = 17.736429841
$ws.Range("R11").Value = 159.627868569
$ws.Range("S11").Value = 0.2029037262835452
$ws.Range("T11").Value = 0.2029037262835451

# Row 12: M2 -> M2
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Icosl"
$ws.Range("C12").Value = "Cd28"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3.0
$ws.Range("F12").Value = 1.0
$ws.Range("G12").Value = 7.390795
$ws.Range("H12").Value = 22.172385
$ws.Range("I12").Value = 0.5493172104644152
$ws.Range("J12").Value = 0.5493172104644152
$ws.Range("K12").Value = 3.0
$ws.Range("L12").Value = 1.0
$ws.Range("M12").Value = 4.065646999999999
$ws.Range("N12").Value = 12.196941
$ws.Range("O12").Value = 0.6257796502061076
$ws.Range("P12").Value = 0.6257796502061076
$ws.Range("Q12").Value = 30.04836351936499
$ws.Range("R12").Value = 270.4352716742849
$ws.Range("S12").Value = 0.3437515318166165
$ws.Range("T12").Value = 0.3437515318166165

# Row 13: M2 -> sCs
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Icosl"
$ws.Range("C13").Value = "Cd28"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3.0
$ws.Range("F13").Value = 1.0
$ws.Range("G13").Value = 7.390795
$ws.Range("H13").Value = 22.172385
$ws.Range("I13").Value = 0.5493172104644152
$ws.Range("J13").Value = 0.5493172104644152
$ws.Range("K13").Value = 1.0
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.03148366666666667
$ws.Range("N13").Value = 0.094451
$ws.Range("O13").Value = 0.004845929298306607
$ws.Range("P13").Value = 0.004845929298306606
$ws.Range("Q13").Value = 0.2326893261816667
$ws.Range("R13").Value = 2.094203935635
$ws.Range("S13").Value = 0.002661952364253566
$ws.Range("T13").Value = 0.002661952364253566

# Row 14: sCs -> M1
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Icosl"
$ws.Range("C14").Value = "Cd28"
$ws.Range("D14").Value = "M1"
$ws.Range("E14").Value = 3.0
$ws.Range("F14").Value = 1.0
$ws.Range("G14").Value = 0.412658
$ws.Range("H14").Value = 1.237974
$ws.Range("I14").Value = 0.03067060328906763
$ws.Range("J14").Value = 0.03067060328906763
$ws.Range("K14").Value = 3.0
$ws.Range("L14").Value = 1.0
$ws.Range("M14").Value = 2.3998
$ws.Range("N14").Value = 7.199400000000001
$ws.Range("O14").Value = 0.3693744204955859
$ws.Range("P14").Value = 0.3693744204955859
$ws.Range("Q14").Value = 0.9902966684000001
$ws.Range("R14").Value = 8.9126700156
$ws.Range("S14").Value = 0.01132893631614937
$ws.Range("T14").Value = 0.01132893631614937

# Row 15: sCs -> M2
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Icosl"
$ws.Range("C15").Value = "Cd28"
$ws.Range("D15").Value = "M2"
$ws.Range("E15").Value = 3.0
$ws.Range("F15").Value = 1.0
$ws.Range("G15").Value = 0.412658
$ws.Range("H15").Value = 1.237974
$ws.Range("I15").Value = 0.03067060328906763
$ws.Range("J15").Value = 0.03067060328906763
$ws.Range("K15").Value = 3.0
$ws.Range("L15").Value = 1.0
$ws.Range("M15").Value = 4.065646999999999
$ws.Range("N15").Value = 12.196941
$ws.Range("O15").Value = 0.6257796502061076
$ws.Range("P15").Value = 0.6257796502061076
$ws.Range("Q15").Value = 1.677721759726
$ws.Range("R15").Value = 15.099495837534
$ws.Range("S15").Value = 0.01919303939784303
$ws.Range("T15").Value = 0.01919303939784303

# Row 16: sCs -> sCs
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Icosl"
$ws.Range("C16").Value = "Cd28"
$ws.Range("D16").Value = "sCs"
$ws.Range("E16").Value = 3.0
$ws.Range("F16").Value = 1.0
$ws.Range("G16").Value = 0.412658
$ws.Range("H16").Value = 1.237974
$ws.Range("I16").Value = 0.03067060328906763
$ws.Range("J16").Value = 0.03067060328906763
$ws.Range("K16").Value = 1.0
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.03148366666666667
$ws.Range("N16").Value = 0.094451
$ws.Range("O16").Value = 0.004845929298306607
$ws.Range("P16").Value = 0.004845929298306606
$ws.Range("Q16").Value = 0.01299198691933333
$ws.Range("R16").Value = 0.116927882274
$ws.Range("S16").Value = 0.0001486275750752318
$ws.Range("T16").Value = 0.0001486275750752318

